# Auto-generated Excel COM-interop script
# Applies scheduled market-price refresh values to the Leve profit tables
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 20000
$ws.Range("J21").Value = 20000
$ws.Range("L21").Value = 20000
$ws.Range("N21").Value = -20936

$ws.Range("H23").Value = 20000
$ws.Range("J23").Value = 20000
$ws.Range("L23").Value = 20000
$ws.Range("N23").Value = -20468

$ws.Range("H58").Value = 3396.6667
$ws.Range("J58").Value = 3396.6667
$ws.Range("L58").Value = 10190.0001
$ws.Range("N58").Value = -10490.0001

$ws.Range("H74").Value = 4164.2856
$ws.Range("I74").Value = 1975
$ws.Range("K74").Value = 1975
$ws.Range("M74").Value = -1039

$ws.Range("H77").Value = 4164.2856
$ws.Range("I77").Value = 1975
$ws.Range("K77").Value = 9875
$ws.Range("M77").Value = -5195

$ws.Range("H98").Value = 1995.25
$ws.Range("I98").Value = 743.625
$ws.Range("J98").Value = 3246.875
$ws.Range("K98").Value = 743.625
$ws.Range("L98").Value = 3246.875
$ws.Range("M98").Value = 754.375
$ws.Range("N98").Value = -6242.875

$ws.Range("H122").Value = 1995.25
$ws.Range("I122").Value = 743.625
$ws.Range("J122").Value = 3246.875
$ws.Range("K122").Value = 2230.875
$ws.Range("L122").Value = 9740.625
$ws.Range("M122").Value = 219.125
$ws.Range("N122").Value = -14640.625

$ws.Range("H132").Value = 41671436
$ws.Range("I132").Value = 55560750
$ws.Range("J132").Value = 3497
$ws.Range("K132").Value = 166682250
$ws.Range("L132").Value = 10491
$ws.Range("M132").Value = -166679720
$ws.Range("N132").Value = -15551

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6392.6665
$ws.Range("I32").Value = 6392.6665
$ws.Range("K32").Value = 6392.6665
$ws.Range("M32").Value = -6105.6665

$ws.Range("H61").Value = 703
$ws.Range("I61").Value = 703
$ws.Range("K61").Value = 703
$ws.Range("M61").Value = -491

$ws.Range("H74").Value = 799.6667
$ws.Range("I74").Value = 799.6667
$ws.Range("K74").Value = 799.6667
$ws.Range("M74").Value = 74.33330000000001

$ws.Range("H77").Value = 799.6667
$ws.Range("I77").Value = 799.6667
$ws.Range("K77").Value = 3998.3335
$ws.Range("M77").Value = 369.6665000000003

$ws.Range("H97").Value = 1018.3077
$ws.Range("I97").Value = 723.2222
$ws.Range("J97").Value = 1682.25
$ws.Range("K97").Value = 723.2222
$ws.Range("L97").Value = 1682.25
$ws.Range("M97").Value = -227.2222
$ws.Range("N97").Value = -2674.25

$ws.Range("H132").Value = 5545.857
$ws.Range("I132").Value = 5545.857
$ws.Range("K132").Value = 16637.571
$ws.Range("M132").Value = -14107.571

$ws.Range("H136").Value = 703
$ws.Range("I136").Value = 703
$ws.Range("K136").Value = 2109
$ws.Range("M136").Value = 441

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7305.647
$ws.Range("I134").Value = 7368.625
$ws.Range("K134").Value = 22105.875
$ws.Range("M134").Value = -19570.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2776.2856
$ws.Range("I31").Value = 1660.7
$ws.Range("K31").Value = 1660.7
$ws.Range("M31").Value = -1365.7

$ws.Range("H34").Value = 2776.2856
$ws.Range("I34").Value = 1660.7
$ws.Range("K34").Value = 1660.7
$ws.Range("M34").Value = -1458.7

$ws.Range("H62").Value = 11066.333
$ws.Range("I62").Value = 11199.625
$ws.Range("K62").Value = 11199.625
$ws.Range("M62").Value = -10575.625

$ws.Range("H65").Value = 11066.333
$ws.Range("I65").Value = 11199.625
$ws.Range("K65").Value = 55998.125
$ws.Range("M65").Value = -52878.125

$ws.Range("H86").Value = 15000
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()

$ws.Range("H89").Value = 15000
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

$ws.Range("H132").Value = 1999.5
$ws.Range("I132").Value = 1999.5
$ws.Range("K132").Value = 5998.5
$ws.Range("M132").Value = -3468.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 1330.6666
$ws.Range("I97").Value = 1522.5
$ws.Range("K97").Value = 4567.5
$ws.Range("M97").Value = -4071.5

$ws.Range("H130").Value = 2349.5
$ws.Range("I130").Value = 1199
$ws.Range("K130").Value = 3597
$ws.Range("M130").Value = 1423

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7710.3335
$ws.Range("I70").Value = 6278.6
$ws.Range("J70").Value = 9500
$ws.Range("K70").Value = 6278.6
$ws.Range("L70").Value = 9500
$ws.Range("M70").Value = -6008.6
$ws.Range("N70").Value = -10040

$ws.Range("H73").Value = 7710.3335
$ws.Range("I73").Value = 6278.6
$ws.Range("J73").Value = 9500
$ws.Range("K73").Value = 6278.6
$ws.Range("L73").Value = 9500
$ws.Range("M73").Value = -5342.6
$ws.Range("N73").Value = -11372

$ws.Range("H132").Value = 3000
$ws.Range("I132").Value = 3000
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 9000
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -6470
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 389
$ws.Range("I93").Value = 389
$ws.Range("K93").Value = 389
$ws.Range("M93").Value = 859

$ws.Range("H132").Value = 15931.4
$ws.Range("I132").Value = 19064.375
$ws.Range("J132").Value = 3399.5
$ws.Range("K132").Value = 57193.125
$ws.Range("L132").Value = 10198.5
$ws.Range("M132").Value = -54663.125
$ws.Range("N132").Value = -15258.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 28886.334
$ws.Range("J94").Value = 28886.334
$ws.Range("L94").Value = 28886.334
$ws.Range("N94").Value = -30688.334

$ws.Range("H100").Value = 574.75
$ws.Range("I100").Value = 433
$ws.Range("K100").Value = 866
$ws.Range("M100").Value = -325

$ws.Range("H132").Value = 900
$ws.Range("I132").Value = 900
$ws.Range("K132").Value = 2700
$ws.Range("M132").Value = -170

$ws.Range("H136").Value = 1268.7368
$ws.Range("I136").Value = 1319.7778
$ws.Range("K136").Value = 3959.3334
$ws.Range("M136").Value = -1409.3334

Write-Output "Applied scheduled market price updates."
